# SQL-Normalization_Relations.xlsx — "Add files via upload"
#
# Renames Sheet1/Sheet2, adds a handful of new label cells that describe
# the MySQL tables referenced by the normalization/relationship diagrams,
# repoints one relabeled cell (Subject -> Subjectid), and updates the
# saved cursor/selection on both sheets.

$wb = $excel.ActiveWorkbook

# --- Rename the first two worksheets -------------------------------------
$wsNorm = $wb.Worksheets.Item("Sheet1")
$wsRel  = $wb.Worksheets.Item("Sheet2")

$wsNorm.Name = "Normalization"
$wsRel.Name  = "Relationships"

# --- Normalization sheet: new/updated text labels -------------------------
# New shared-string text must be written in this exact order so the
# resulting sharedStrings.xml index assignment lines up with the workbook
# being reproduced (Subjectid, subject table, student_subject table,
# address table, student table).
$wsNorm.Range("K20").Value = "Subjectid"
$wsNorm.Range("K19").Value = "subject table"
$wsNorm.Range("A30").Value = "student_subject table"
$wsNorm.Range("H28").Value = "address table"
$wsNorm.Range("B20").Value = "student table"

# --- Restore the saved selections on each sheet ---------------------------
$wsNorm.Activate()
$wsNorm.Range("I1").Select()

$wsRel.Activate()
$wsRel.Range("D28").Select()
